# Update LR-pairs TPM-derived values for Sertad1 -> Ar sheet.
# The underlying per-cell TPM values were recomputed, which changes the
# per-cluster ligand/receptor average (and therefore derived) statistics
# for the "ECs" and "MuSCs" clusters (the "FAPs" cluster stats are
# unchanged). This updates columns G:T for rows 2-10 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> column letter -> new value
$updates = @{
    2  = @{ G = 10.38120433333333; H = 31.143613; I = 0.2980959031155602; J = 0.2980959031155602; K = 3; L = 1; M = 0.4797913333333333; N = 1.439374; O = 0.0238859545710277; P = 0.0238859545710277; Q = 4.980811868695778; R = 44.827306818262; S = 0.007120305199627745; T = 0.007120305199627745 }
    3  = @{ G = 10.38120433333333; H = 31.143613; I = 0.2980959031155602; J = 0.2980959031155602; O = 0.7934729695107392; P = 0.7934729695107391; Q = 165.4587248031566; R = 1489.128523228409; S = 0.2365310414440892; T = 0.2365310414440891 }
    4  = @{ G = 10.38120433333333; H = 31.143613; I = 0.2980959031155602; J = 0.2980959031155602; M = 3.668666666666667; N = 11.006; O = 0.1826410759182331; P = 0.1826410759182331; Q = 38.08517829755556; R = 342.766604678; S = 0.0544445564718433; T = 0.0544445564718433 }
    5  = @{ I = 0.4170467778049336; J = 0.4170467778049336; K = 3; L = 1; M = 0.4797913333333333; N = 1.439374; O = 0.0238859545710277; P = 0.0238859545710277; Q = 6.96833307328911; R = 62.71499765960199; S = 0.009961560388642128; T = 0.009961560388642128 }
    6  = @{ I = 0.4170467778049336; J = 0.4170467778049336; O = 0.7934729695107392; P = 0.7934729695107391; S = 0.3309153452097661; T = 0.3309153452097661 }
    7  = @{ I = 0.4170467778049336; J = 0.4170467778049336; M = 3.668666666666667; N = 11.006; O = 0.1826410759182331; P = 0.1826410759182331; Q = 53.28251990422222; R = 479.542679138; S = 0.07616987220652538; T = 0.07616987220652538 }
    8  = @{ G = 9.920170000000001; H = 29.76051; I = 0.2848573190795062; J = 0.2848573190795062; K = 3; L = 1; M = 0.4797913333333333; N = 1.439374; O = 0.0238859545710277; P = 0.0238859545710277; Q = 4.759611591193333; R = 42.83650432074; S = 0.006804088982757828; T = 0.006804088982757828 }
    9  = @{ G = 9.920170000000001; H = 29.76051; I = 0.2848573190795062; J = 0.2848573190795062; O = 0.7934729695107392; P = 0.7934729695107391; Q = 158.1106223639367; R = 1422.99560127543; S = 0.226026582856884; T = 0.2260265828568839 }
    10 = @{ G = 9.920170000000001; H = 29.76051; I = 0.2848573190795062; J = 0.2848573190795062; M = 3.668666666666667; N = 11.006; O = 0.1826410759182331; P = 0.1826410759182331; Q = 36.39379700666667; R = 327.54417306; S = 0.05202664723986446; T = 0.05202664723986446 }
}

foreach ($r in $updates.Keys) {
    $cols = $updates[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
